$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create the new "Assumptions" worksheet. Worksheets.Add() inserts the new
# sheet before the currently active sheet, which places it before "Hazards"
# (matching the target sheet order: Assumptions, Hazards).
# ---------------------------------------------------------------------------
$assumptions = $wb.Worksheets.Add()
$assumptions.Name = "Assumptions"

# Fill in the assumption rows. The order in which the cells are written
# controls the order new entries are appended to the shared string table,
# so it is intentionally not a simple top-to-bottom, left-to-right sweep.
$assumptions.Range("A1").Value = "A-1"
$assumptions.Range("A3").Value = "A-2"
$assumptions.Range("A2").Value = "A-1.1"
$assumptions.Range("A4").Value = "A-3"
$assumptions.Range("B3").Value = "The vehicle has an SAE automation level of 3."
$assumptions.Range("B1").Value = "The vehicle is used exclusively in urban areas."
$assumptions.Range("B4").Value = "The vehicle has a normal braking distance of (speed [km/h] / 10)^2."
$assumptions.Range("B5").Value = "The vehicle has an emergency braking distance of (speed [km/h] / 10)^2 / 2."
$assumptions.Range("A5").Value = "A-4"
$assumptions.Range("A6").Value = "A-5"
$assumptions.Range("B6").Value = "The vehicle has a minimum (front) sensing range of 100 m."
$assumptions.Range("B2").Value = "The vehicle drives at a maximum speed of 50 km/h (~ 13.89 m/s)."

# Column B is wide enough to show the full assumption text.
$assumptions.Columns.Item(2).ColumnWidth = 68.14

# Match the page setup used by the rest of the workbook.
$assumptions.PageSetup.PaperSize = 9
$assumptions.PageSetup.Orientation = 1
$assumptions.PageSetup.TopMargin = 56.6929133866
$assumptions.PageSetup.BottomMargin = 56.6929133866

# Give the descriptive cells a touched font style (mirrors the extra cellXfs
# entry present in the edited workbook).
$assumptions.Range("B1:B2").Font.ThemeColor = 1

# Leave the selection on B7, just below the data, as in the authored file.
$assumptions.Range("B7").Select()

# ---------------------------------------------------------------------------
# Fix the "accerates" -> "accelerates" typo in the Hazards sheet.
# ---------------------------------------------------------------------------
$hazards = $wb.Worksheets.Item("Hazards")
$hazards.Range("C4").Value = "While driving, the vehicle suddenly accelerates without any reason and can collide with an obstacle on or offside the road."

# Make Hazards the active tab again (it is now the second sheet).
$hazards.Activate()
$hazards.Range("A2").Select()
